# Generate Report for Archive
#
# Refresh the localization-status report:
#  - The "zh-cn" / "de-de" status cells move from "Ready for handoff" to
#    "In Translation" (Overview!E2:F2, zh-cn!C2, de-de!C2 all share the
#    same shared-string value).
#  - The corresponding Status columns narrow to fit the new, shorter text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Update every cell that currently shows "Ready for handoff"
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# Narrow the Status columns to match the new (shorter) text
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
